# Updates the "LOS Galacticos" roster sheet:
#  - Adds "Dennis Schröder" and "Amen Thompson" back into the list in new
#    positions (row 2 and row 9), shifting the rows that followed them down.
#  - The net effect on the data grid (Player | Position | Team) is written
#    here explicitly, row by row, so the final sheet content matches
#    regardless of how rows got reshuffled along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dennis Schröder",        "PG,SG",    "Golden State Warriors"),
    @("Anfernee Simons",        "PG,SG",    "Portland Trail Blazers"),
    @("Fred VanVleet",          "PG",       "Houston Rockets"),
    @("Dillon Brooks",          "SG,SF",    "Houston Rockets"),
    @("Anthony Edwards",        "SG,SF",    "Minnesota Timberwolves"),
    @("Giannis Antetokounmpo",  "PF,C",     "Milwaukee Bucks"),
    @("Bobby Portis",           "PF,C",     "Milwaukee Bucks"),
    @("Amen Thompson",          "PF,C",     "Houston Rockets"),
    @("Jaren Jackson Jr.",      "PF,C",     "Memphis Grizzlies"),
    @("Ivica Zubac",            "C",        "LA Clippers"),
    @("Andrew Wiggins",         "SF,PF",    "Golden State Warriors"),
    @("Bilal Coulibaly",        "SG,SF",    "Washington Wizards"),
    @("James Harden",           "PG,SG",    "LA Clippers"),
    @("Zion Williamson",        "PF,C",     "New Orleans Pelicans"),
    @("Jayson Tatum",           "SF,PF",    "Boston Celtics"),
    @("Jonathan Kuminga",       "SF,PF",    "Golden State Warriors"),
    @("Jerami Grant",           "SF,PF",    "Portland Trail Blazers"),
    @("Paul George",            "SG,SF,PF", "Philadelphia 76ers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Range("A$row").Value = $entry[0]
    $ws.Range("B$row").Value = $entry[1]
    $ws.Range("C$row").Value = $entry[2]
    $row++
}
